$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" footer timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 18:35"

# --- Simple data refreshes (no reordering of the country list) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1437053
$ws.Range("C4").Value = 6705
$ws.Range("D4").Value = 310834
$ws.Range("E4").Value = 1040731
$ws.Range("F4").Value = 16337
$ws.Range("G4").Value = 291
$ws.Range("H4").Value = 85488

# Row 8: Italia
$ws.Range("B8").Value = 223096
$ws.Range("C8").Value = 992
$ws.Range("D8").Value = 115288
$ws.Range("E8").Value = 76440
$ws.Range("F8").Value = 855
$ws.Range("G8").Value = 262
$ws.Range("H8").Value = 31368

# Row 66: Luxemburgo
$ws.Range("B66").Value = 3915
$ws.Range("C66").Value = 11
$ws.Range("D66").Value = 3665
$ws.Range("E66").Value = 147
$ws.Range("F66").Value = 11
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 103

# --- Countries that swapped position in the ranking: Moldavia now
#     ranks above Ghana, Santa Lucia now ranks above Belice, and
#     Sahara Occidental now ranks above San Bartolome. Update both the
#     country names and their data rows to reflect the new order. ---

# Rows 62/63: Moldavia now above Ghana, with Moldavia's figures updated
$ws.Range("A62").Value = "Moldavia"
$ws.Range("B62").Value = 5553
$ws.Range("C62").Value = 147
$ws.Range("D62").Value = 2228
$ws.Range("E62").Value = 3131
$ws.Range("F62").Value = 251
$ws.Range("G62").Value = 9
$ws.Range("H62").Value = 194

$ws.Range("A63").Value = "Ghana"
$ws.Range("B63").Value = 5408
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 514
$ws.Range("E63").Value = 4870
$ws.Range("F63").Value = 5
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 24

# Rows 194/195: Santa Lucia now above Belice
$ws.Range("A194").Value = "Santa Lucia"
$ws.Range("B194").Value = 18
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 18
$ws.Range("E194").Value = 0
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0

$ws.Range("A195").Value = "Belice"
$ws.Range("B195").Value = 18
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 16
$ws.Range("E195").Value = 0
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 2

# Rows 215/216: Sahara Occidental now above San Bartolome (figures unchanged)
$ws.Range("A215").Value = "Sahara Occidental"
$ws.Range("A216").Value = "San Bartolome"
